$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 114.052635
$ws.Range("I33").Value = 121.411766
$ws.Range("J33").Value = 51.5
$ws.Range("K33").Value = 121.411766
$ws.Range("L33").Value = 51.5
$ws.Range("M33").Value = 107.588234
$ws.Range("N33").Value = -509.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 35000
$ws.Range("J87").Value = 35000
$ws.Range("L87").Value = 35000
$ws.Range("N87").Value = -37496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H90").Value = 35000
$ws.Range("J90").Value = 35000
$ws.Range("L90").Value = 105000
$ws.Range("N90").Value = -117480

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 1510.1538
$ws.Range("I100").Value = 1508.2
$ws.Range("J100").Value = 1516.6666
$ws.Range("K100").Value = 1508.2
$ws.Range("L100").Value = 1516.6666
$ws.Range("M100").Value = -967.2
$ws.Range("N100").Value = -2598.6666

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H106").Value = 3010
$ws.Range("I106").Value = 2900
$ws.Range("J106").Value = 4000
$ws.Range("K106").Value = 2900
$ws.Range("L106").Value = 4000
$ws.Range("M106").Value = -2269
$ws.Range("N106").Value = -5262

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 45324.168
$ws.Range("J133").Value = 45324.168
$ws.Range("L133").Value = 45324.168
$ws.Range("N133").Value = -55444.168

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 15320.286
$ws.Range("I28").Value = 4888.4
$ws.Range("J28").Value = 41400
$ws.Range("K28").Value = 4888.4
$ws.Range("L28").Value = 41400
$ws.Range("M28").Value = -4696.4
$ws.Range("N28").Value = -41784

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5397.0547
$ws.Range("I32").Value = 4629.755
$ws.Range("J32").Value = 11663.333
$ws.Range("K32").Value = 4629.755
$ws.Range("L32").Value = 11663.333
$ws.Range("M32").Value = -4342.755
$ws.Range("N32").Value = -12237.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 15320.286
$ws.Range("I99").Value = 4888.4
$ws.Range("J99").Value = 41400
$ws.Range("K99").Value = 4888.4
$ws.Range("L99").Value = 41400
$ws.Range("M99").Value = -1893.4
$ws.Range("N99").Value = -47390

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 158.125
$ws.Range("I22").Value = 200
$ws.Range("J22").Value = 116.25
$ws.Range("K22").Value = 200
$ws.Range("L22").Value = 116.25
$ws.Range("M22").Value = -27
$ws.Range("N22").Value = -462.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 160755.44
$ws.Range("J132").Value = 160755.44
$ws.Range("L132").Value = 160755.44
$ws.Range("N132").Value = -170875.44

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 85568.336
$ws.Range("I134").Value = 2611
$ws.Range("J134").Value = 251483
$ws.Range("K134").Value = 7833
$ws.Range("L134").Value = 754449
$ws.Range("M134").Value = -5298
$ws.Range("N134").Value = -759519

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H135").Value = 70779
$ws.Range("J135").Value = 70779
$ws.Range("L135").Value = 70779
$ws.Range("N135").Value = -80919

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 49940
$ws.Range("J138").Value = 49940
$ws.Range("L138").Value = 49940
$ws.Range("N138").Value = -60220

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4195.1177
$ws.Range("I134").Value = 2751.8
$ws.Range("J134").Value = 6257
$ws.Range("K134").Value = 8255.400000000001
$ws.Range("L134").Value = 18771
$ws.Range("M134").Value = -5720.400000000001
$ws.Range("N134").Value = -23841

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 55522.332
$ws.Range("J140").Value = 55522.332
$ws.Range("L140").Value = 55522.332
$ws.Range("N140").Value = -65882.33199999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H136").Value = 3672.2222
$ws.Range("I136").Value = 2610
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 7830
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -2730
$ws.Range("N136").Value = -25200

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1014854.94
$ws.Range("I122").Value = 2195188.2
$ws.Range("J122").Value = 3140.5715
$ws.Range("K122").Value = 6585564.600000001
$ws.Range("L122").Value = 9421.7145
$ws.Range("M122").Value = -6583114.600000001
$ws.Range("N122").Value = -14321.7145

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2722.6
$ws.Range("J126").Value = 1664.25
$ws.Range("L126").Value = 4992.75
$ws.Range("N126").Value = -9932.75

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H133").Value = 35976
$ws.Range("J133").Value = 35976
$ws.Range("L133").Value = 35976
$ws.Range("N133").Value = -46096

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H136").Value = 22249.5
$ws.Range("J136").Value = 22249.5
$ws.Range("L136").Value = 66748.5
$ws.Range("N136").Value = -71848.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 61925
$ws.Range("J138").Value = 73850
$ws.Range("L138").Value = 73850
$ws.Range("N138").Value = -84130

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H139").Value = 64708.668
$ws.Range("J139").Value = 64708.668
$ws.Range("L139").Value = 64708.668
$ws.Range("N139").Value = -74988.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 2379.4473
$ws.Range("I136").Value = 1279.7778
$ws.Range("J136").Value = 5078.636
$ws.Range("K136").Value = 3839.3334
$ws.Range("L136").Value = 15235.908
$ws.Range("M136").Value = -1289.3334
$ws.Range("N136").Value = -20335.908

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1672.0454
$ws.Range("I81").Value = 1496.5385
$ws.Range("J81").Value = 1925.5555
$ws.Range("K81").Value = 2993.077
$ws.Range("L81").Value = 3851.111
$ws.Range("M81").Value = -1932.077
$ws.Range("N81").Value = -5973.111

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1672.0454
$ws.Range("I84").Value = 1496.5385
$ws.Range("J84").Value = 1925.5555
$ws.Range("K84").Value = 14965.385
$ws.Range("L84").Value = 19255.555
$ws.Range("M84").Value = -9661.385000000002
$ws.Range("N84").Value = -29863.555

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H133").Value = 75783
$ws.Range("J133").Value = 75783
$ws.Range("L133").Value = 75783
$ws.Range("N133").Value = -85903

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 1601.3334
$ws.Range("I136").Value = 1592.4615
$ws.Range("J136").Value = 1624.4
$ws.Range("K136").Value = 4777.3845
$ws.Range("L136").Value = 4873.200000000001
$ws.Range("M136").Value = -2227.3845
$ws.Range("N136").Value = -9973.200000000001

